$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New formula cells in column H (rows 12-15), matching the style (s="0")
# of their row-mates (e.g. D12/D13/D14 which already use the default style).
$ws.Range("H12").Formula = "=5/6"
$ws.Range("H13").Formula = "=72/6"
$ws.Range("H14").Formula = "=144/6"
$ws.Range("H15").Formula = "=1/12"

# Move/restore the active selection from H12 to H16.
$ws.Range("H16").Select() | Out-Null
